$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column E (Obrigatorio) to "S" for rows 2 through 10
for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 5).Value = "S"
}
